$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.00977443609022556
$ws.Range("C2").Value = 0.0075187969924812
$ws.Range("D2").Value = 0.0120300751879699
$ws.Range("E2").Value = 0.968421052631579
$ws.Range("F2").Value = 0.00150375939849624
$ws.Range("G2").Value = 0.0037593984962406
$ws.Range("H2").Value = 0.00075187969924812
$ws.Range("I2").Value = 0.0195488721804511
$ws.Range("J2").Value = 0.00150375939849624
$ws.Range("K2").Value = 0.00225563909774436
$ws.Range("L2").Value = 0.00075187969924812
$ws.Range("M2").Value = 0.00601503759398496
$ws.Range("N2").Value = 0.00075187969924812
$ws.Range("P2").Value = 0.00075187969924812
$ws.Range("Q2").Value = 0.857894736842105
$ws.Range("R2").Value = 0.639097744360902
$ws.Range("S2").Value = 0.00075187969924812
$ws.Range("U2").Value = 0.946616541353383
$ws.Range("V2").Value = 0.0142857142857143
$ws.Range("W2").Value = 0.0330827067669173
$ws.Range("X2").Value = 0.0225563909774436

# Row 3
$ws.Range("B3").Value = 0.00300751879699248
$ws.Range("C3").Value = 0.00075187969924812
$ws.Range("D3").Value = 0.00300751879699248
$ws.Range("E3").Value = 0.0075187969924812
$ws.Range("G3").Value = 0.954135338345865
$ws.Range("H3").Value = 0.0037593984962406
$ws.Range("I3").Value = 0.97218045112782
$ws.Range("J3").Value = 0.0172932330827068
$ws.Range("K3").Value = 0.978195488721805
$ws.Range("L3").Value = 0.0172932330827068
$ws.Range("M3").Value = 0.0075187969924812
$ws.Range("N3").Value = 0.00150375939849624
$ws.Range("O3").Value = 0.970676691729323
$ws.Range("P3").Value = 0.00300751879699248
$ws.Range("R3").Value = 0.00225563909774436
$ws.Range("S3").Value = 0.989473684210526
$ws.Range("T3").Value = 0.0037593984962406
$ws.Range("U3").Value = 0.00150375939849624
$ws.Range("V3").Value = 0.00451127819548872
$ws.Range("W3").Value = 0.00225563909774436
$ws.Range("X3").Value = 0.0075187969924812

# Row 4
$ws.Range("B4").Value = 0.968421052631579
$ws.Range("C4").Value = 0.969924812030075
$ws.Range("D4").Value = 0.983458646616541
$ws.Range("E4").Value = 0.0075187969924812
$ws.Range("F4").Value = 0.990977443609023
$ws.Range("G4").Value = 0.0075187969924812
$ws.Range("H4").Value = 0.0218045112781955
$ws.Range("I4").Value = 0.0075187969924812
$ws.Range("J4").Value = 0.00150375939849624
$ws.Range("K4").Value = 0.00075187969924812
$ws.Range("L4").Value = 0.0037593984962406
$ws.Range("M4").Value = 0.969924812030075
$ws.Range("N4").Value = 0.988721804511278
$ws.Range("O4").Value = 0.0142857142857143
$ws.Range("Q4").Value = 0.13984962406015
$ws.Range("R4").Value = 0.354135338345865
$ws.Range("S4").Value = 0.00075187969924812
$ws.Range("T4").Value = 0.00150375939849624
$ws.Range("U4").Value = 0.0503759398496241
$ws.Range("V4").Value = 0.979699248120301
$ws.Range("W4").Value = 0.96390977443609
$ws.Range("X4").Value = 0.969924812030075

# Row 5
$ws.Range("B5").Value = 0.0180451127819549
$ws.Range("C5").Value = 0.0210526315789474
$ws.Range("D5").Value = 0.00150375939849624
$ws.Range("E5").Value = 0.0165413533834586
$ws.Range("F5").Value = 0.0075187969924812
$ws.Range("G5").Value = 0.0345864661654135
$ws.Range("H5").Value = 0.973684210526316
$ws.Range("I5").Value = 0.00075187969924812
$ws.Range("J5").Value = 0.979699248120301
$ws.Range("K5").Value = 0.018796992481203
$ws.Range("L5").Value = 0.978195488721805
$ws.Range("M5").Value = 0.0157894736842105
$ws.Range("N5").Value = 0.00902255639097744
$ws.Range("O5").Value = 0.0150375939849624
$ws.Range("P5").Value = 0.996240601503759
$ws.Range("Q5").Value = 0.00225563909774436
$ws.Range("R5").Value = 0.00451127819548872
$ws.Range("S5").Value = 0.00902255639097744
$ws.Range("T5").Value = 0.994736842105263
$ws.Range("U5").Value = 0.00150375939849624
$ws.Range("V5").Value = 0.00075187969924812
$ws.Range("W5").Value = 0.00075187969924812

